# Applies the cryptos.xlsx price/coin-list update described in the commit diff
# ("Updated symbol list on Wed Dec 14 19:26:51 UTC 2022 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # The D column stores numeric-looking quotes as literal text (inlineStr).
    # Briefly force a text format so Excel does not coerce the assignment to a
    # number, then restore the original "General" format so no visible style
    # change is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "270.90"
Set-TextValue $ws.Range("D3") "23.02"
Set-TextValue $ws.Range("D4") "6.378"
Set-TextValue $ws.Range("D5") "0.06240"
Set-TextValue $ws.Range("D6") "3.665"
Set-TextValue $ws.Range("D7") "6.756"
Set-TextValue $ws.Range("D8") "1.395"
Set-TextValue $ws.Range("D9") "0.8383"
Set-TextValue $ws.Range("D10") "0.01371"
Set-TextValue $ws.Range("D11") "0.1620"
Set-TextValue $ws.Range("D12") "0.08359"
Set-TextValue $ws.Range("D13") "0.03423"
Set-TextValue $ws.Range("D14") "0.03183"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09316"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D16") "3.939"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D17") "0.001733"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D18") "0.04888"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D19") "0.006265"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D20") "0.005425"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D21") "0.001097"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D22") "0.0001511"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D23") "3.738"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D24") "2.317"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D25") "0.3337"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D26") "0.1251"
$ws.Range("E26").Value = "25ProBitTokenPROB"
Set-TextValue $ws.Range("D27") "0.0002700"
Set-TextValue $ws.Range("D40") "0.04686"
Set-TextValue $ws.Range("D41") "0.006934"
Set-TextValue $ws.Range("D42") "0.1167"
Set-TextValue $ws.Range("D43") "0.003484"
Set-TextValue $ws.Range("D44") "0.01254"
Set-TextValue $ws.Range("D45") "0.00006291"
Set-TextValue $ws.Range("D46") "0.00000000755"
Set-TextValue $ws.Range("D47") "0.7046"
Set-TextValue $ws.Range("D48") "0.1288"
Set-TextValue $ws.Range("D50") "0.01248"

Write-Host "Applied cryptos.xlsx price update."
